$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price values in column D are stored as text (e.g. "43.764.19", "0.619"),
# not numbers, so force text formatting before assigning to avoid Excel
# auto-converting numeric-looking strings into numbers.

$dCells = @("D2","D3","D5","D6","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D20","D21","D24","D26","D27","D28","D29","D32","D34","D35","D39","D40","D41","D43","D44","D45","D50")
foreach ($ref in $dCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = '43.764.19'
$ws.Range("E2").Value = '  +0.09%  '
$ws.Range("D3").Value = '2.287.77'
$ws.Range("E3").Value = '  -0.18%  '
$ws.Range("E4").Value = '  +0.32%  '
$ws.Range("D5").Value = '115.58'
$ws.Range("E5").Value = '  +12.00%  '
$ws.Range("D6").Value = '268.50'
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("E7").Value = '  +0.32%  '
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("D9").Value = '0.619'
$ws.Range("E9").Value = '  +1.72%  '
$ws.Range("D10").Value = '49.49'
$ws.Range("E10").Value = '  +8.99%  '
$ws.Range("D11").Value = '0.0944'
$ws.Range("E11").Value = '  +0.98%  '
$ws.Range("D12").Value = '8.94'
$ws.Range("E12").Value = '  +11.65%  '
$ws.Range("D13").Value = '0.107'
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("D14").Value = '15.89'
$ws.Range("E14").Value = '  +0.52%  '
$ws.Range("D15").Value = '2.631.23'
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("D16").Value = '0.884'
$ws.Range("E16").Value = '  +3.14%  '
$ws.Range("D17").Value = '2.281.17'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '43.609.32'
$ws.Range("E18").Value = '  -0.24%  '
$ws.Range("E19").Value = '  -1.03%  '
$ws.Range("D20").Value = '7.04'
$ws.Range("E20").Value = '  +12.38%  '
$ws.Range("D21").Value = '72.27'
$ws.Range("E22").Value = '  -3.49%  '
$ws.Range("E23").Value = '  +8.34%  '
$ws.Range("D24").Value = '233.14'
$ws.Range("E24").Value = '  -0.14%  '
$ws.Range("E25").Value = '  +1.35%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("D27").Value = '11.67'
$ws.Range("E27").Value = '  +4.08%  '
$ws.Range("D28").Value = '3.92'
$ws.Range("E28").Value = '  +0.89%  '
$ws.Range("D29").Value = '41.99'
$ws.Range("E29").Value = '  +6.54%  '
$ws.Range("E30").Value = '  -1.71%  '
$ws.Range("E31").Value = '  +0.36%  '
$ws.Range("D32").Value = '173.14'
$ws.Range("E32").Value = '  -2.36%  '
$ws.Range("E33").Value = '  +3.76%  '
$ws.Range("D34").Value = '21.64'
$ws.Range("E34").Value = '  -0.67%  '
$ws.Range("D35").Value = '5.70'
$ws.Range("E35").Value = '  +4.36%  '
$ws.Range("E36").Value = '  +0.43%  '
$ws.Range("E37").Value = '  -3.76%  '
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("D39").Value = '0.108'
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("D40").Value = '3.81'
$ws.Range("E40").Value = '  +6.73%  '
$ws.Range("D41").Value = '14.59'
$ws.Range("E41").Value = '  +19.17%  '
$ws.Range("E42").Value = '  +4.50%  '
$ws.Range("D43").Value = '74.47'
$ws.Range("E43").Value = '  +13.33%  '
$ws.Range("D44").Value = '0.242'
$ws.Range("E44").Value = '  +2.60%  '
$ws.Range("D45").Value = '6.39'
$ws.Range("E45").Value = '  +22.08%  '
$ws.Range("E46").Value = '  +0.20%  '
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("E49").Value = '  +3.97%  '
$ws.Range("D50").Value = '102.81'
$ws.Range("E50").Value = '  +4.21%  '
$ws.Range("E51").Value = '  -1.46%  '
